$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (also reflected in the workbook's sheet tab name / tag)
$ws.Name = "Through 2022-04-27"

# Update the header label for April in column A, row 5
$ws.Range("A5").Value = "April (through 04-27)"

# Update April row (row 5) values
$ws.Range("B5").Value = 20
$ws.Range("C5").Value = 32
$ws.Range("D5").Value = 55
$ws.Range("F5").Value = 40
$ws.Range("G5").Value = 55
$ws.Range("H5").Value = 93
$ws.Range("I5").Value = 109

# Update Total row (row 6) values
$ws.Range("B6").Value = 86
$ws.Range("C6").Value = 160
$ws.Range("D6").Value = 244
$ws.Range("F6").Value = 150
$ws.Range("G6").Value = 253
$ws.Range("H6").Value = 516
$ws.Range("I6").Value = 544
